$d = $word.ActiveDocument

# --- 1. Insert a space between "${cargo_representante}" and "de la ${" but
#        only in the FIRST occurrence (the "Vista la solicitud..." paragraph).
#        The second occurrence, in the later "APRUEBESE..." paragraph, is
#        left untouched, so the change must be scoped precisely.
$old = '$' + '{cargo_representante}de la $' + '{'
$new = '$' + '{cargo_representante} de la $' + '{'

$fullText = $d.Content.Text
$firstIdx = $fullText.IndexOf($old)

$editRange = $d.Range($firstIdx, $firstIdx + $old.Length)
$editRange.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null

# --- 2. Relocate the hidden "_GoBack" bookmark. In the source document it
#        sits between the "anio_resolucion" run and the closing "}" run
#        inside the paragraph just edited above; after the edit it must
#        instead wrap the "${nombre_delegado}" run near the end of the
#        document (in the paragraph right after "ANTE MI:").
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

$fullText = $d.Content.Text
$needle = '$' + '{nombre_delegado}'
$start = $fullText.IndexOf($needle)
$end = $start + $needle.Length

$target = $d.Range($start, $end)
$d.Bookmarks.Add("_GoBack", $target) | Out-Null
